$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 990.6896400000001
$ws.Range("I19").Value = 496.3846
$ws.Range("J19").Value = 1392.3125
$ws.Range("K19").Value = 496.3846
$ws.Range("L19").Value = 1392.3125
$ws.Range("M19").Value = -321.3846
$ws.Range("N19").Value = -1742.3125

$ws.Range("H41").Value = 13333973
$ws.Range("I41").Value = 397
$ws.Range("J41").Value = 22223024
$ws.Range("K41").Value = 397
$ws.Range("L41").Value = 22223024
$ws.Range("M41").Value = 43
$ws.Range("N41").Value = -22223904

$ws.Range("H98").Value = 1072.9348
$ws.Range("I98").Value = 1121.5
$ws.Range("J98").Value = 563
$ws.Range("K98").Value = 1121.5
$ws.Range("L98").Value = 563
$ws.Range("M98").Value = 376.5
$ws.Range("N98").Value = -3559

$ws.Range("H103").Value = 1633.6666
$ws.Range("J103").Value = 1000
$ws.Range("L103").Value = 3000
$ws.Range("N103").Value = -4172

$ws.Range("H111").Value = 15882875
$ws.Range("I111").Value = 22230868
$ws.Range("J111").Value = 12890
$ws.Range("K111").Value = 66692604
$ws.Range("L111").Value = 38670
$ws.Range("M111").Value = -66689537
$ws.Range("N111").Value = -44804

$ws.Range("H112").Value = 6320.7915
$ws.Range("J112").Value = 7411.316
$ws.Range("L112").Value = 22233.948
$ws.Range("N112").Value = -24449.948

$ws.Range("H122").Value = 1072.9348
$ws.Range("I122").Value = 1121.5
$ws.Range("J122").Value = 563
$ws.Range("K122").Value = 3364.5
$ws.Range("L122").Value = 1689
$ws.Range("M122").Value = -914.5
$ws.Range("N122").Value = -6589

$ws.Range("H129").Value = 142858180
$ws.Range("I129").Value = 166667380
$ws.Range("K129").Value = 500002140
$ws.Range("M129").Value = -499997140

$ws.Range("H131").Value = 4132.129
$ws.Range("I131").Value = 1959
$ws.Range("J131").Value = 9444.223
$ws.Range("K131").Value = 5877
$ws.Range("L131").Value = 28332.669
$ws.Range("M131").Value = -837
$ws.Range("N131").Value = -38412.669

$ws.Range("H138").Value = 4910.485
$ws.Range("J138").Value = 5315.2
$ws.Range("L138").Value = 15945.6
$ws.Range("N138").Value = -26225.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 3540.6
$ws.Range("I3").Value = 3540.6
$ws.Range("K3").Value = 3540.6
$ws.Range("M3").Value = -3425.6

$ws.Range("H32").Value = 6472.95
$ws.Range("I32").Value = 5208.231
$ws.Range("K32").Value = 5208.231
$ws.Range("M32").Value = -4921.231

$ws.Range("H61").Value = 6110
$ws.Range("I61").Value = 7165.8335
$ws.Range("J61").Value = 3998.3333
$ws.Range("K61").Value = 7165.8335
$ws.Range("L61").Value = 3998.3333
$ws.Range("M61").Value = -6953.8335
$ws.Range("N61").Value = -4422.3333

$ws.Range("H74").Value = 47846.074
$ws.Range("I74").Value = 4494.2295
$ws.Range("K74").Value = 4494.2295
$ws.Range("M74").Value = -3620.2295

$ws.Range("H77").Value = 47846.074
$ws.Range("I77").Value = 4494.2295
$ws.Range("K77").Value = 22471.1475
$ws.Range("M77").Value = -18103.1475

$ws.Range("H136").Value = 6110
$ws.Range("I136").Value = 7165.8335
$ws.Range("J136").Value = 3998.3333
$ws.Range("K136").Value = 21497.5005
$ws.Range("L136").Value = 11994.9999
$ws.Range("M136").Value = -18947.5005
$ws.Range("N136").Value = -17094.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3761959.2
$ws.Range("I107").Value = 4763934
$ws.Range("J107").Value = 4553.25
$ws.Range("K107").Value = 4763934
$ws.Range("L107").Value = 4553.25
$ws.Range("M107").Value = -4762014
$ws.Range("N107").Value = -8393.25

$ws.Range("H109").Value = 59564.668
$ws.Range("J109").Value = 59564.668
$ws.Range("L109").Value = 59564.668
$ws.Range("N109").Value = -62338.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 264.6
$ws.Range("I7").Value = 102.181816
$ws.Range("K7").Value = 102.181816
$ws.Range("M7").Value = 10.818184

$ws.Range("H31").Value = 3313.66
$ws.Range("I31").Value = 1209.2632
$ws.Range("J31").Value = 3807.284
$ws.Range("K31").Value = 1209.2632
$ws.Range("L31").Value = 3807.284
$ws.Range("M31").Value = -914.2632000000001
$ws.Range("N31").Value = -4397.284

$ws.Range("H34").Value = 3313.66
$ws.Range("I34").Value = 1209.2632
$ws.Range("J34").Value = 3807.284
$ws.Range("K34").Value = 1209.2632
$ws.Range("L34").Value = 3807.284
$ws.Range("M34").Value = -1007.2632
$ws.Range("N34").Value = -4211.284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5243.091
$ws.Range("J5").Value = 10640.5
$ws.Range("L5").Value = 31921.5
$ws.Range("N5").Value = -32145.5

$ws.Range("H56").Value = 8965.083000000001
$ws.Range("I56").Value = 8965.083000000001
$ws.Range("K56").Value = 8965.083000000001
$ws.Range("M56").Value = -8435.083000000001

$ws.Range("H68").Value = 2606.8276
$ws.Range("I68").Value = 2423.6667
$ws.Range("J68").Value = 2689.25
$ws.Range("K68").Value = 7271.000100000001
$ws.Range("L68").Value = 8067.75
$ws.Range("M68").Value = -6460.000100000001
$ws.Range("N68").Value = -9689.75

$ws.Range("H71").Value = 2606.8276
$ws.Range("I71").Value = 2423.6667
$ws.Range("J71").Value = 2689.25
$ws.Range("K71").Value = 21813.0003
$ws.Range("L71").Value = 24203.25
$ws.Range("M71").Value = -17757.0003
$ws.Range("N71").Value = -32315.25

$ws.Range("H107").Value = 1278.3529
$ws.Range("I107").Value = 815.6667
$ws.Range("J107").Value = 1530.7273
$ws.Range("K107").Value = 2447.0001
$ws.Range("L107").Value = 4592.1819
$ws.Range("M107").Value = -527.0001000000002
$ws.Range("N107").Value = -8432.1819

$ws.Range("H132").Value = 1263.5652
$ws.Range("I132").Value = 835.6
$ws.Range("K132").Value = 7520.400000000001
$ws.Range("M132").Value = -4990.400000000001

$ws.Range("H135").Value = 5243.091
$ws.Range("J135").Value = 10640.5
$ws.Range("L135").Value = 95764.5
$ws.Range("N135").Value = -100834.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 10000.5
$ws.Range("I34").Value = 10001
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 10001
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -9733
$ws.Range("N34").Value = -10536

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H76").Value = 10000.5
$ws.Range("I76").Value = 10001
$ws.Range("J76").Value = 10000
$ws.Range("K76").Value = 10001
$ws.Range("L76").Value = 10000
$ws.Range("M76").Value = -9686
$ws.Range("N76").Value = -10630

$ws.Range("H79").Value = 10000.5
$ws.Range("I79").Value = 10001
$ws.Range("J79").Value = 10000
$ws.Range("K79").Value = 10001
$ws.Range("L79").Value = 10000
$ws.Range("M79").Value = -8909
$ws.Range("N79").Value = -12184

$ws.Range("H80").Value = 104879430
$ws.Range("I80").Value = 131098664
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 131098664
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -131097666
$ws.Range("N80").Value = -4496

$ws.Range("H83").Value = 104879430
$ws.Range("I83").Value = 131098664
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 655493320
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -655488328
$ws.Range("N83").Value = -22484

$ws.Range("H126").Value = 5155939
$ws.Range("I126").Value = 2529730
$ws.Range("K126").Value = 7589190
$ws.Range("M126").Value = -7586720

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 55556932
$ws.Range("I82").Value = 55556932
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 55556932
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -55556571
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 55556932
$ws.Range("I85").Value = 55556932
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 55556932
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -55555684
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 125008330
$ws.Range("I107").Value = 142862370
$ws.Range("K107").Value = 428587110
$ws.Range("M107").Value = -428585190

$ws.Range("H109").Value = 34992.5
$ws.Range("J109").Value = 34992.5
$ws.Range("L109").Value = 34992.5
$ws.Range("N109").Value = -37766.5

$ws.Range("H126").Value = 2084.9412
$ws.Range("I126").Value = 2165.25
$ws.Range("J126").Value = 800
$ws.Range("K126").Value = 6495.75
$ws.Range("L126").Value = 2400
$ws.Range("M126").Value = -4025.75
$ws.Range("N126").Value = -7340

$ws.Range("H135").Value = 35714
$ws.Range("J135").Value = 40428
$ws.Range("L135").Value = 40428
$ws.Range("N135").Value = -50568

$ws.Range("H136").Value = 5224.054
$ws.Range("I136").Value = 6886.684
$ws.Range("J136").Value = 3469.0557
$ws.Range("K136").Value = 20660.052
$ws.Range("L136").Value = 10407.1671
$ws.Range("M136").Value = -18110.052
$ws.Range("N136").Value = -15507.1671
